$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = 27.24147881842065
    3  = 30.02570492936247
    4  = 32.84022063337129
    5  = 29.86484597673077
    6  = 22.68941707443092
    7  = 23.82025847019073
    8  = 22.3872337490217
    9  = 21.7378231240628
    10 = 20.68101164237834
    11 = 18.61287903799967
    12 = 10.94230523340259
    13 = 38.74181640808195
    14 = 40.3418180623393
    15 = 41.69329990042433
    16 = 35.01762949195202
    18 = 31.71476036536174
    19 = 28.50321598913881
    20 = 27.55183100771479
    21 = 25.80972703483493
    22 = 24.87392167116274
    23 = 15.21328379750986
    34 = 7.836074300910229
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 4).Value = $values[$row]
}
